# Insert a new data row at row 24 (pushes the existing rows 24-109 down to 25-110)
# and populate it with the new "1a (guarda)" price observation for Asterix potatoes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(24).Insert()

$ws.Range("A24").Value = 11
$ws.Range("B24").Value = "Vega Monumental Concepción"
$ws.Range("C24").Value = "Bíobío"
$ws.Range("D24").Value = 44453
$ws.Range("E24").Value = 8
$ws.Range("F24").Value = 100114001
$ws.Range("G24").Value = "Papa"
$ws.Range("H24").Value = "Asterix"
$ws.Range("I24").Value = "1a (guarda)"
$ws.Range("J24").Value = 2000
$ws.Range("K24").Value = 8500
$ws.Range("L24").Value = 9000
$ws.Range("M24").Value = 8750
$ws.Range("N24").Value = "`$/saco 25 kilos"
$ws.Range("O24").Value = "Provincia de Arauco"
$ws.Range("P24").Value = 350
$ws.Range("Q24").Value = 25
$ws.Range("R24").Value = "Hortaliza"
